$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values for rows 2-25 (columns B, C, D, E, F, I, M),
# recomputed for the 380 kV case ("case with 380 kV done").
$data = @{
    2 = @(13.65895741983901, 11.2043421740263, 5.586399925725559, 10.32407531646839, 33.62818623086564, 25.14875172362442, 15.77829257748993)
    3 = @(13.15024845339819, 10.59599414025422, 5.591683893182625, 10.20885581676488, 33.05243794692844, 24.95110829641694, 15.52038966883093)
    4 = @(12.83248428438939, 10.20746421422591, 5.595648450364507, 10.14062121320593, 32.70431194916925, 24.83528284550419, 15.36518962880919)
    5 = @(12.7018695763316, 10.04556421834095, 5.597443460257481, 10.11347197480437, 32.56397969273642, 24.78950267544478, 15.30282391434682)
    6 = @(12.68012065065981, 10.01847178260564, 5.597752307453222, 10.10900425618849, 32.54077528743601, 24.78198745753725, 15.29252375947475)
    7 = @(12.83072698592323, 10.20529493959321, 5.595671934285869, 10.14025237663655, 32.70241294423975, 24.8346596535733, 15.36434486837613)
    8 = @(13.48481472922662, 10.99780735550255, 5.588071362942139, 10.28384239341818, 33.42865767258559, 25.07947543259222, 15.68876553340731)
    9 = @(14.71514471797221, 12.42573084990086, 5.578955306243857, 10.58414530038612, 34.88767319538274, 25.60199335889818, 16.34612047903157)
    10 = @(15.57613804973007, 13.39029481927168, 5.575887657324521, 10.81450480607957, 35.97033394873438, 26.00966560446753, 16.83682367501076)
    11 = @(15.9565892874352, 13.80960601628261, 5.575301518490003, 10.92103699651739, 36.46301176795999, 26.19979741469472, 17.06066679499266)
    12 = @(16.09890918209372, 13.96551694292151, 5.575197635071243, 10.96159599427753, 36.64942479311842, 26.27242406880279, 17.14543301364174)
    13 = @(16.0683377941323, 13.93206752090488, 5.575214730177832, 10.95285173363883, 36.60928699625619, 26.25675542708172, 17.12717837611481)
    14 = @(15.96833379208801, 13.82249084809289, 5.575290595958884, 10.92436963692556, 36.47835214910703, 26.20576013297164, 17.0676410189105)
    15 = @(15.90684699436932, 13.75499590023677, 5.57535249246757, 10.90695089177028, 36.39812563306383, 26.17460444790128, 17.03117034529852)
    16 = @(15.55103693038724, 13.36249341152685, 5.575942394237687, 10.80757490680988, 35.93812437175373, 25.99733039093944, 16.82219992993927)
    17 = @(15.32978350985084, 13.11666140962261, 5.576512785499131, 10.74703351810175, 35.65584356837437, 25.88974547933243, 16.69409895210836)
    18 = @(15.20147533831053, 12.97343566591451, 5.576916939895633, 10.71237751125715, 35.4935093965065, 25.82830875601758, 16.62048205523952)
    19 = @(15.15785663744832, 12.92463011156065, 5.577066794536111, 10.70067305321808, 35.438555868739, 25.80758481858025, 16.59557023982356)
    20 = @(15.35344592453053, 13.14302058600216, 5.576444180007934, 10.75346132932972, 35.68589143010167, 25.90115252174445, 16.70772960140127)
    21 = @(15.99775584501317, 13.85475465561813, 5.575265094484147, 10.93272987549718, 36.51681643518535, 26.22072201295261, 17.08512921973647)
    22 = @(16.40860581160182, 14.30314344750958, 5.575183488531004, 11.05114521370705, 37.05890226934631, 26.43321700941573, 17.33175323279508)
    23 = @(16.1903050849514, 14.06538474865292, 5.575163448426378, 10.98784088834978, 36.76972679851896, 26.31948679637968, 17.20015602917027)
    24 = @(15.34275158330508, 13.13110949498587, 5.576474959243225, 10.75055484604014, 35.67230692520238, 25.89599410097914, 16.70156708455363)
    25 = @(14.38913099362217, 12.05383093295343, 5.580791205327126, 10.5010712863937, 34.49033375985838, 25.45630182426903, 16.1665582578723)
}

$columns = @("B", "C", "D", "E", "F", "I", "M")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}
